$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (reordered)
$ws.Range("A1").Value = "Variable Name"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Unit"
$ws.Range("F1").Value = "Fixed Costs"

# Row 2
$ws.Range("A2").Value = "container_cost"
$ws.Range("B2").Value = 2.0065108174152999
$ws.Range("C2").Value = "$"
$ws.Range("F2").Value = "Test Equipment"

# Row 3
$ws.Range("A3").Value = "cost_lid"
$ws.Range("F3").Value = "Manufacturing Tooling"

# Row 4
$ws.Range("A4").Value = "Battery_Cost_USD"
$ws.Range("F4").Value = "QA Equipment"

# Row 5
$ws.Range("A5").Value = "heating element"
$ws.Range("F5").Value = "safety and compliance testing "

# Row 6 (new data in A/B), F6 label stays text but index shifts automatically
$ws.Range("A6").Value = "unit quantity"
$ws.Range("B6").Value = 0
$ws.Range("F6").Value = "Total_fixed_costs"

# Row 8
$ws.Range("F8").Value = "Fixed_costs_per_unit"
$ws.Range("A8").Value = "Total_Variable_Cost"
$ws.Range("B8").Formula = "=SUM(B2:B5)"
$ws.Range("G8").Formula = "=G6 / B6"

# Column width
$ws.Columns.Item(1).ColumnWidth = 19.33

# Selection
$ws.Range("C15").Select()
